# Apply "branding" vs "configuration" terminology edits to the document.
#
# We locate each target passage with Find.Execute (no in-place replacement
# string, so Word's AutoCorrect / AutoFormat "smart quotes" machinery never
# gets a chance to touch straight apostrophes) and then assign the new text
# directly to the found Range's .Text property, which writes the literal
# characters we supply.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
    if ($ok) {
        $r.Text = $new
    } else {
        Write-Output "NOT FOUND: $old"
    }
}

# Longer, full-sentence replacements first (these consume every occurrence
# of "...Configuration folder..." etc. that lives inside a full sentence).

Replace-Text `
    "The Profiles Open Source branding system aims to help you to use configurations in order to switch between different variations for the overall appearance (along with some smaller details), i.e., the 'branding', of the pages of Open Source Profiles." `
    "The Profiles Open Source branding system aims to help you to use 'branding' configurations in order to switch between different variations for the overall appearance (along with some smaller details), i.e., the branding, of the pages of Open Source Profiles."

Replace-Text `
    "Profiles contains two Branding Configurations." `
    "Profiles contains two options for branding."

Replace-Text `
    "The default configuration, OpenSource branding, is intended for sites wanting to get up and running quickly with a simple and functional interface. " `
    "The default branding, OpenSource, is intended for sites wanting to get up and running quickly with a simple and functional interface. "

Replace-Text `
    "The Foo configuration represents a fictional (satirical) institution. Foo is a good sandbox for experimenting with the features of the branding setup, and we use it below to illustrate the branding system." `
    "The Foo branding represents a fictional (satirical) institution. Foo is a good sandbox for experimenting with the features of the branding setup, and we use it below to illustrate the branding system."

Replace-Text `
    "The current configuration resides in the Configuration folder, located just above (i.e., containing as a sub-folder) StaticFiles. This Configuration folder could contain several variants for the branding. In this Open Source release, it contains the sub-folders OpenSource and Foo, and the configuration files in-effect are copies of those in OpenSource." `
    "The current branding resides in the Branding folder. This folder is located just above (i.e., containing as a sub-folder) StaticFiles. The Branding folder could contain several variants for the branding. In this Open Source release, it contains the sub-folders OpenSource and Foo, and the branding files in-effect are copies of those in OpenSource."

Replace-Text `
    "If you are not sure which branding configuration is currently in-effect, then in the Configuration folder, you can check by running (in a shell) the script, checkConfigurationFolder. For example, if you are currently using the Foo configuration then " `
    "If you are not sure which branding is currently in-effect, then in the Branding folder, you can check by running (in a shell) the script, checkBrandingFolder. For example, if you are currently using Foo branding then "

Replace-Text `
    "% checkConfigurationFolder.bash   Foo" `
    "% checkBrandingFolder.bash   Foo"

Replace-Text `
    "should indicate three empty diff results. If you are using OpenSource or some other added configuration, then you will likely see many diff results." `
    "should indicate three empty diff results. If you are using OpenSource or some other custom branding, then you will likely see many diff results."

Replace-Text `
    "For your institution's branding, you can provide your own custom configuration. Your configuration will use your own headers and footers, and your own particular values for properties, css definitions and text-snippets. Compose it by modifying / providing these files:" `
    "For your institution's branding, you can provide your own custom files. Your branding can use your own headers and footers, and your own particular values for properties, css definitions and text-snippets. Compose it by modifying / providing these files:"

Replace-Text `
    "The following examples assume that the files of the included 'Foo' branding -- that is, the three 'myBranding' files from the Foo folder -- have been copied up into the Configuration folder." `
    "The following examples assume that the files of the included 'Foo' branding -- that is, the three 'myBranding' files from the Foo folder -- have been copied up into the Branding folder."

Replace-Text `
    "The data and functions from the three myBranding.* configuration files get used across the application, so we shall refer to a few files ('general-application files') beyond those three." `
    "The data and functions from the three myBranding.* files get used across the application, so we shall refer to a few files ('general-application files') beyond those three."

# Finally, the lone heading run "Configuration" (in "Some Examples from The
# Foo Configuration") -- every other standalone occurrence of the word has
# already been consumed by the full-sentence replacements above.

Replace-Text `
    "Configuration" `
    "Branding"

Write-Output "Done applying branding/configuration edits."
